# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-09-28 Saturday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-09-29 Sunday", 2)

# Update the table of division problems. Addressed by (row, column) rather
# than by text search because several old values repeat (e.g. "40÷3=13, 1"
# and "70÷8=8, 6" each occur twice) but map to different replacements
# depending on position.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91÷3=30, 1"
$t.Cell(1,2).Range.Text = "36÷8=4, 4"
$t.Cell(1,3).Range.Text = "83÷6=13, 5"
$t.Cell(1,4).Range.Text = "97÷5=19, 2"
$t.Cell(1,5).Range.Text = "31÷5=6, 1"

$t.Cell(5,1).Range.Text = "29÷7=4, 1"
$t.Cell(5,2).Range.Text = "80÷3=26, 2"
$t.Cell(5,3).Range.Text = "40÷3=13, 1"
$t.Cell(5,4).Range.Text = "37÷9=4, 1"
$t.Cell(5,5).Range.Text = "77÷2=38, 1"

$t.Cell(9,1).Range.Text = "45÷3=15, 0"
$t.Cell(9,2).Range.Text = "68÷8=8, 4"
$t.Cell(9,3).Range.Text = "23÷8=2, 7"
$t.Cell(9,4).Range.Text = "67÷6=11, 1"
$t.Cell(9,5).Range.Text = "82÷6=13, 4"

$t.Cell(13,1).Range.Text = "26÷5=5, 1"
$t.Cell(13,2).Range.Text = "13÷6=2, 1"
$t.Cell(13,3).Range.Text = "23÷5=4, 3"
$t.Cell(13,4).Range.Text = "99÷6=16, 3"
$t.Cell(13,5).Range.Text = "82÷6=13, 4"

$t.Cell(17,1).Range.Text = "63÷6=10, 3"
$t.Cell(17,2).Range.Text = "43÷5=8, 3"
$t.Cell(17,3).Range.Text = "50÷9=5, 5"
$t.Cell(17,4).Range.Text = "96÷7=13, 5"
$t.Cell(17,5).Range.Text = "17÷2=8, 1"
